$d = $word.ActiveDocument
$p1 = $d.Paragraphs(1)

# --- "last minute" paragraph formatting tweaks on the ID placeholder paragraph ---
# Add a thin paragraph border (5 twips of space on every edge) and widen the
# left indent from 120 -> 225 twips (6pt -> 11.25pt).
$pf = $p1.Range.ParagraphFormat
$pf.Borders.DistanceFromTop = 5
$pf.Borders.DistanceFromLeft = 5
$pf.Borders.DistanceFromBottom = 5
$pf.Borders.DistanceFromRight = 5
$pf.LeftIndent = 11.25

# --- locate the old placeholder token's extent without disturbing formatting ---
$findRange = $d.Content
$findRange.Find.ClearFormatting()
$findRange.Find.Execute("**ID__AFFARS_5327_topic_5__ID**", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0)
$idStart = $findRange.Start
$idEnd = $findRange.End

# --- drop the trailing space run that sits right after the placeholder ---
$spaceRange = $d.Range($idEnd, $idEnd + 1)
$spaceRange.Delete()

# --- rewrite the placeholder id (topic_5 -> 90) in place, in its own run ---
$idRange = $d.Range($idStart, $idEnd)
$idRange.Text = "**ID__AFFARS_5327_90__ID**"
